$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New frequency values for rows 2-5, columns B-X (updated run for publication)
$data = @(
    @(0.0118606375092661,0.00593031875463306,0.00593031875463306,0.00889547813194959,0.00148257968865827,0.011119347664937,0.116382505559674,0.0429948109710897,0.0148257968865827,0.0022238695329874,0.00889547813194959,0.00148257968865827,0.00148257968865827,0.000741289844329133,0.0185322461082283,0.00148257968865827,0.00815418828762046,0.982209043736101,0.00815418828762046,0.0422535211267606,0.0289103039288362,0.00815418828762046,0.0022238695329874),
    @(0.00148257968865827,0.0022238695329874,0.687916975537435,0.0022238695329874,0.0378057820607858,0.00593031875463306,0.00518902891030393,0.916234247590808,0.97405485544848,0.00148257968865827,0.00518902891030393,0.00370644922164566,0.00518902891030393,0.983691623424759,0.00148257968865827,0.000741289844329133,0,0,0.981467753891772,0.0896960711638251,0,0.0022238695329874,0.00296515937731653),
    @(0.98295033358043,0.99110452186805,0.00963676797627873,0.985174203113417,0.0192735359525575,0.979243884358784,0.865826538176427,0.00296515937731653,0,0.022238695329874,0.98295033358043,0.994069681245367,0.988880652335063,0.0022238695329874,0.972572275759822,0.997776130467013,0.985915492957746,0.0170496664195701,0.00148257968865827,0.0259451445515196,0.966641957005189,0.987398072646405,0.992587101556709),
    @(0.00370644922164566,0.000741289844329133,0.296515937731653,0.00370644922164566,0.941438102297999,0.00370644922164566,0.0118606375092661,0.0370644922164566,0.011119347664937,0.97405485544848,0.0022238695329874,0.000741289844329133,0.00296515937731653,0.0133432171979244,0.00667160859896219,0,0.00593031875463306,0.000741289844329133,0.00889547813194959,0.842105263157895,0.0044477390659748,0.0022238695329874,0.00148257968865827)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $colNum = $j + 2   # Column B = 2
        $ws.Cells.Item($rowNum, $colNum).Value = $rowValues[$j]
    }
}
